$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 1).Value = "2025-10-13 01:50:37"
}
